$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = '@'
    $rng.Value = $value
}

Set-TextValue 'D2' '42.616.73'
Set-TextValue 'E2' '  -0.48%  '
Set-TextValue 'D3' '2.285.56'
Set-TextValue 'E3' '  -2.38%  '
Set-TextValue 'E4' '  -0.03%  '
Set-TextValue 'D5' '300.38'
Set-TextValue 'E5' '  -2.17%  '
Set-TextValue 'D6' '96.81'
Set-TextValue 'E6' '  -4.54%  '
Set-TextValue 'E7' '  -1.28%  '
Set-TextValue 'E8' '  +0.02%  '
Set-TextValue 'D9' '0.498'
Set-TextValue 'E9' '  -2.98%  '
Set-TextValue 'D10' '33.17'
Set-TextValue 'E10' '  -5.50%  '
Set-TextValue 'D11' '0.0787'
Set-TextValue 'E11' '  -1.14%  '
Set-TextValue 'D12' '49.36'
Set-TextValue 'E12' '  -5.92%  '
Set-TextValue 'E13' '  +0.20%  '
Set-TextValue 'E14' '  -2.71%  '
Set-TextValue 'D15' '2.639.94'
Set-TextValue 'E15' '  -2.67%  '
Set-TextValue 'D16' '15.41'
Set-TextValue 'E16' '  -0.05%  '
Set-TextValue 'D17' '2.329.68'
Set-TextValue 'E17' '  -0.83%  '
Set-TextValue 'D18' '0.786'
Set-TextValue 'E18' '  -1.50%  '
Set-TextValue 'D19' '42.551.93'
Set-TextValue 'E19' '  -0.54%  '
Set-TextValue 'E20' '  -1.01%  '
Set-TextValue 'D21' '11.47'
Set-TextValue 'E21' '  -2.05%  '
Set-TextValue 'D22' '6.00'
Set-TextValue 'E22' '  -3.83%  '
Set-TextValue 'D23' '66.75'
Set-TextValue 'E23' '  -0.99%  '
Set-TextValue 'D24' '234.33'
Set-TextValue 'E24' '  -1.25%  '
Set-TextValue 'D25' '1.93'
Set-TextValue 'E25' '  -3.07%  '
Set-TextValue 'E26' '  -3.44%  '
Set-TextValue 'E27' '  +0.08%  '
Set-TextValue 'D28' '24.37'
Set-TextValue 'E28' '  -3.76%  '
Set-TextValue 'D29' '165.99'
Set-TextValue 'E29' '  +3.69%  '
Set-TextValue 'E30' '  -5.85%  '
Set-TextValue 'D31' '33.86'
Set-TextValue 'E31' '  -4.26%  '
Set-TextValue 'D32' '9.10'
Set-TextValue 'E32' '  -2.84%  '
Set-TextValue 'E33' '  +0.00%  '
Set-TextValue 'D34' '4.96'
Set-TextValue 'E34' '  -3.64%  '
Set-TextValue 'E35' '  -3.92%  '
Set-TextValue 'E36' '  -4.45%  '
Set-TextValue 'D37' '4.32'
Set-TextValue 'E37' '  -5.49%  '
Set-TextValue 'E38' '  -6.22%  '
Set-TextValue 'D39' '16.19'
Set-TextValue 'E39' '  -9.22%  '
Set-TextValue 'D40' '1.76'
Set-TextValue 'E40' '  -6.35%  '
Set-TextValue 'D41' '0.0994'
Set-TextValue 'E41' '  -3.83%  '
Set-TextValue 'E42' '  -2.42%  '
Set-TextValue 'D43' '2.45'
Set-TextValue 'E43' '  -2.26%  '
Set-TextValue 'D44' '1.961.48'
Set-TextValue 'E44' '  -3.00%  '
Set-TextValue 'E45' '  -1.08%  '
Set-TextValue 'D46' '17.93'
Set-TextValue 'E46' '  -6.37%  '
Set-TextValue 'D47' '9.70'
Set-TextValue 'E47' '  -7.54%  '
Set-TextValue 'E48' '  -6.06%  '
Set-TextValue 'B49' 'HuobiToken'
Set-TextValue 'C49' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D49' '2.82'
Set-TextValue 'E49' '  -3.48%  '
Set-TextValue 'B50' 'MultiversX'
Set-TextValue 'C50' 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue 'D50' '53.11'
Set-TextValue 'E50' '  -6.78%  '
Set-TextValue 'D51' '2.509.87'
Set-TextValue 'E51' '  -2.98%  '
